$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.06056433333333333
$ws.Range("H2").Value = 0.181693
$ws.Range("M2").Value = 91.60947133333333
$ws.Range("N2").Value = 274.828414
$ws.Range("O2").Value = 0.7121576185577153
$ws.Range("P2").Value = 0.7121576185577152
$ws.Range("Q2").Value = 5.548266558322444
$ws.Range("R2").Value = 49.934399024902
$ws.Range("S2").Value = 0.7121576185577153
$ws.Range("T2").Value = 0.7121576185577152

# Row 3
$ws.Range("G3").Value = 0.06056433333333333
$ws.Range("H3").Value = 0.181693
$ws.Range("M3").Value = 22.83185066666667
$ws.Range("N3").Value = 68.495552
$ws.Range("O3").Value = 0.1774912152792038
$ws.Range("P3").Value = 0.1774912152792038
$ws.Range("Q3").Value = 1.382795814392889
$ws.Range("R3").Value = 12.445162329536
$ws.Range("S3").Value = 0.1774912152792038
$ws.Range("T3").Value = 0.1774912152792038

# Row 4
$ws.Range("G4").Value = 0.06056433333333333
$ws.Range("H4").Value = 0.181693
$ws.Range("M4").Value = 7.077809999999999
$ws.Range("N4").Value = 21.23343
$ws.Range("O4").Value = 0.05502178149094856
$ws.Range("P4").Value = 0.05502178149094855
$ws.Range("Q4").Value = 0.4286628441099999
$ws.Range("R4").Value = 3.85796559699
$ws.Range("S4").Value = 0.05502178149094856
$ws.Range("T4").Value = 0.05502178149094855

# Row 5
$ws.Range("G5").Value = 0.06056433333333333
$ws.Range("H5").Value = 0.181693
$ws.Range("M5").Value = 7.117379
$ws.Range("N5").Value = 21.352137
$ws.Range("O5").Value = 0.05532938467213248
$ws.Range("P5").Value = 0.05532938467213247
$ws.Range("Q5").Value = 0.4310593142156666
$ws.Range("R5").Value = 3.879533827941
$ws.Range("S5").Value = 0.05532938467213248
$ws.Range("T5").Value = 0.05532938467213247
